$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column sometimes holds numeric-looking text (e.g. "1.00", "69.516.82")
# that must stay literal text instead of being auto-converted to a number by Excel,
# so format the cells we are about to rewrite as Text first.
$ws.Range("D2:D12").NumberFormat = "@"
$ws.Range("D14:D19").NumberFormat = "@"
$ws.Range("D21:D36").NumberFormat = "@"
$ws.Range("D39:D44").NumberFormat = "@"
$ws.Range("D46:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.516.82"
$ws.Range("E2").Value = "  +6.05%  "

$ws.Range("D3").Value = "3.571.46"
$ws.Range("E3").Value = "  +5.22%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "592.00"
$ws.Range("E5").Value = "  +5.54%  "

$ws.Range("D6").Value = "192.53"
$ws.Range("E6").Value = "  +9.38%  "

$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("D8").Value = "3.562.21"
$ws.Range("E8").Value = "  +5.32%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "0.184"
$ws.Range("E10").Value = "  +6.22%  "

$ws.Range("D11").Value = "0.662"
$ws.Range("E11").Value = "  +4.02%  "

$ws.Range("D12").Value = "58.30"
$ws.Range("E12").Value = "  +9.12%  "

$ws.Range("E13").Value = "  +5.82%  "

$ws.Range("D14").Value = "9.70"
$ws.Range("E14").Value = "  +5.17%  "

$ws.Range("D15").Value = "4.128.19"
$ws.Range("E15").Value = "  +4.90%  "

$ws.Range("D16").Value = "19.29"
$ws.Range("E16").Value = "  +5.41%  "

$ws.Range("D17").Value = "3.559.15"
$ws.Range("E17").Value = "  +4.65%  "

$ws.Range("D18").Value = "69.364.06"
$ws.Range("E18").Value = "  +6.04%  "

$ws.Range("D19").Value = "12.46"
$ws.Range("E19").Value = "  +5.20%  "

$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").Value = "1.05"
$ws.Range("E21").Value = "  +4.67%  "

$ws.Range("D22").Value = "502.75"
$ws.Range("E22").Value = "  +4.38%  "

$ws.Range("D23").Value = "5.55"
$ws.Range("E23").Value = "  +12.53%  "

$ws.Range("D24").Value = "17.21"
$ws.Range("E24").Value = "  +19.89%  "

$ws.Range("D25").Value = "4.43"
$ws.Range("E25").Value = "  +7.80%  "

$ws.Range("D26").Value = "91.30"
$ws.Range("E26").Value = "  +1.95%  "

$ws.Range("D27").Value = "3.05"
$ws.Range("E27").Value = "  +4.62%  "

$ws.Range("D28").Value = "11.21"
$ws.Range("E28").Value = "  +5.35%  "

$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  +6.73%  "

$ws.Range("D30").Value = "32.08"
$ws.Range("E30").Value = "  +2.52%  "

$ws.Range("D31").Value = "7.47"
$ws.Range("E31").Value = "  +13.94%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "12.17"
$ws.Range("E32").Value = "  +5.72%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "618.70"
$ws.Range("E33").Value = "  +7.47%  "

$ws.Range("D34").Value = "65.36"

$ws.Range("D35").Value = "0.115"
$ws.Range("E35").Value = "  +6.49%  "

$ws.Range("D36").Value = "0.0₃0832"
$ws.Range("E36").Value = "  +12.42%  "

$ws.Range("E37").Value = "  +4.90%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").Value = "38.00"
$ws.Range("E39").Value = "  +6.00%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.398"
$ws.Range("E40").Value = "  +6.37%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "3.64"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").Value = "3.331.23"
$ws.Range("E42").Value = "  +7.52%  "

$ws.Range("D43").Value = "3.07"
$ws.Range("E43").Value = "  +9.77%  "

$ws.Range("D44").Value = "2.73"
$ws.Range("E44").Value = "  +12.16%  "

$ws.Range("E45").Value = "  +5.86%  "

$ws.Range("D46").Value = "2.93"
$ws.Range("E46").Value = "  +23.32%  "

$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  +3.76%  "

$ws.Range("D48").Value = "0.138"
$ws.Range("E48").Value = "  +2.27%  "

$ws.Range("D49").Value = "9.08"
$ws.Range("E49").Value = "  +7.62%  "

$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "3.23"
$ws.Range("E50").Value = "  +4.14%  "

$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.10%  "
